$d = $word.ActiveDocument

# --- Part 1: split "Low-cost archival storage." run and insert _GoBack bookmark ---
$findRange = $d.Content
[void]$findRange.Find.Execute("Low-cost archival stora", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($findRange.End, $findRange.End)
[void]$d.Bookmarks.Add("_GoBack", $bmRange)

# --- Part 2: replace the trailing empty bookmark paragraph with the new
#     "Additional Database Services" section, ending in a fresh empty paragraph ---
$tailPara = $d.Paragraphs.Last
$tailRange = $tailPara.Range
$tailRange.Collapse(1)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="4096"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading1"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
    <w:t>Additional Database Services</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:rFonts w:ascii="unset" w:hAnsi="unset" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
    <w:t xml:space="preserve">Amazon </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:rFonts w:ascii="unset" w:hAnsi="unset" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
    <w:t>DocumentDB</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:hyperlink r:id="rId5" w:tgtFrame="_blank" w:history="1">
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Strong"/>
        <w:rFonts w:ascii="unset" w:eastAsiaTheme="majorEastAsia" w:hAnsi="unset" w:cs="Arial"/>
        <w:color w:val="0000FF"/>
        <w:sz w:val="21"/>
        <w:szCs w:val="21"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t xml:space="preserve">Amazon </w:t>
    </w:r>
    <w:proofErr w:type="spellStart"/>
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Strong"/>
        <w:rFonts w:ascii="unset" w:eastAsiaTheme="majorEastAsia" w:hAnsi="unset" w:cs="Arial"/>
        <w:color w:val="0000FF"/>
        <w:sz w:val="21"/>
        <w:szCs w:val="21"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>DocumentDB</w:t>
    </w:r>
    <w:proofErr w:type="spellEnd"/>
  </w:hyperlink>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> is a document database service that supports MongoDB workloads. (MongoDB is a document database program.)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="27"/>
      <w:szCs w:val="27"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:rFonts w:ascii="unset" w:hAnsi="unset" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Amazon Neptune</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:hyperlink r:id="rId6" w:tgtFrame="_blank" w:history="1">
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Strong"/>
        <w:rFonts w:ascii="unset" w:eastAsiaTheme="majorEastAsia" w:hAnsi="unset" w:cs="Arial"/>
        <w:color w:val="0000FF"/>
        <w:sz w:val="21"/>
        <w:szCs w:val="21"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>Amazon Neptune</w:t>
    </w:r>
  </w:hyperlink>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> is a graph database service. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>You can use Amazon Neptune to build and run applications that work with highly connected datasets, such as recommendation engines, fraud detection, and knowledge graphs.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="27"/>
      <w:szCs w:val="27"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:rFonts w:ascii="unset" w:hAnsi="unset" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
    <w:t xml:space="preserve">Amazon Quantum Ledger Database (Amazon QLDB) </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:hyperlink r:id="rId7" w:tgtFrame="_blank" w:history="1">
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Strong"/>
        <w:rFonts w:ascii="unset" w:eastAsiaTheme="majorEastAsia" w:hAnsi="unset" w:cs="Arial"/>
        <w:color w:val="0000FF"/>
        <w:sz w:val="21"/>
        <w:szCs w:val="21"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>Amazon Quantum Ledger Database (Amazon QLDB)</w:t>
    </w:r>
  </w:hyperlink>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> is a ledger database service. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>You can use Amazon QLDB to review a complete history of all the changes that have been made to your application data.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="27"/>
      <w:szCs w:val="27"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:rFonts w:ascii="unset" w:hAnsi="unset" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
    <w:t>Amazon Managed Blockchain</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:hyperlink r:id="rId8" w:tgtFrame="_blank" w:history="1">
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Strong"/>
        <w:rFonts w:ascii="unset" w:eastAsiaTheme="majorEastAsia" w:hAnsi="unset" w:cs="Arial"/>
        <w:color w:val="0000FF"/>
        <w:sz w:val="21"/>
        <w:szCs w:val="21"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>Amazon Managed Blockchain</w:t>
    </w:r>
  </w:hyperlink>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> is a service that you can use to create and manage blockchain networks with open-source frameworks. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>Blockchain is a distributed ledger system that lets multiple parties run transactions and share data without a central authority.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="27"/>
      <w:szCs w:val="27"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:rFonts w:ascii="unset" w:hAnsi="unset" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
    <w:t xml:space="preserve">Amazon </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:rFonts w:ascii="unset" w:hAnsi="unset" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
    <w:t>ElastiCache</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:hyperlink r:id="rId9" w:tgtFrame="_blank" w:history="1">
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Strong"/>
        <w:rFonts w:ascii="unset" w:eastAsiaTheme="majorEastAsia" w:hAnsi="unset" w:cs="Arial"/>
        <w:color w:val="0000FF"/>
        <w:sz w:val="21"/>
        <w:szCs w:val="21"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t xml:space="preserve">Amazon </w:t>
    </w:r>
    <w:proofErr w:type="spellStart"/>
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Strong"/>
        <w:rFonts w:ascii="unset" w:eastAsiaTheme="majorEastAsia" w:hAnsi="unset" w:cs="Arial"/>
        <w:color w:val="0000FF"/>
        <w:sz w:val="21"/>
        <w:szCs w:val="21"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>ElastiCache</w:t>
    </w:r>
    <w:proofErr w:type="spellEnd"/>
  </w:hyperlink>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> is a service that adds caching layers on top of your databases to help improve the read times of common requests. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>It supports two types of data stores: Redis and Memcached.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading3"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="27"/>
      <w:szCs w:val="27"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="Strong"/>
      <w:rFonts w:ascii="unset" w:hAnsi="unset" w:cs="Arial"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
      <w:color w:val="1F1F1F"/>
    </w:rPr>
    <w:t>Amazon DynamoDB Accelerator</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:hyperlink r:id="rId10" w:tgtFrame="_blank" w:history="1">
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Strong"/>
        <w:rFonts w:ascii="unset" w:eastAsiaTheme="majorEastAsia" w:hAnsi="unset" w:cs="Arial"/>
        <w:color w:val="0000FF"/>
        <w:sz w:val="21"/>
        <w:szCs w:val="21"/>
        <w:u w:val="single"/>
      </w:rPr>
      <w:t>Amazon DynamoDB Accelerator (DAX)</w:t>
    </w:r>
  </w:hyperlink>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t xml:space="preserve"> is an in-memory cache for DynamoDB. </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="NormalWeb"/>
    <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    <w:spacing w:before="0" w:beforeAutospacing="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="1F1F1F"/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>It helps improve response times from single-digit milliseconds to microseconds.</w:t>
  </w:r>
</w:p>
<w:p/>

</w:body>
</w:document>
</pkg:xmlData></pkg:part>
<pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
<Relationship Id="rId5" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://aws.amazon.com/documentdb/" TargetMode="External"/>
<Relationship Id="rId6" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://aws.amazon.com/neptune/" TargetMode="External"/>
<Relationship Id="rId7" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://aws.amazon.com/qldb/" TargetMode="External"/>
<Relationship Id="rId8" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://aws.amazon.com/managed-blockchain/" TargetMode="External"/>
<Relationship Id="rId9" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://aws.amazon.com/elasticache/" TargetMode="External"/>
<Relationship Id="rId10" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://aws.amazon.com/dynamodb/dax/" TargetMode="External"/>
</Relationships>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$tailRange.InsertXML($xml)

# Remove the now-stale empty paragraph that used to hold the _GoBack bookmark
$oldTail = $d.Paragraphs.Last.Previous
if ($oldTail.Range.Text -eq "\r") {
    $oldTail.Range.Delete()
}

Write-Output "done"
